$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must stay stored as text (matches original inlineStr cells).
$textCells = @("D5", "D6", "D7", "D8", "D11", "D12", "D13", "D14", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '63.333.68'
$ws.Range("E2").Value = '  +6.21%  '
$ws.Range("D3").Value = '2.433.64'
$ws.Range("E3").Value = '  +5.74%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '566.24'
$ws.Range("E5").Value = '  +4.52%  '
$ws.Range("D6").Value = '141.04'
$ws.Range("E6").Value = '  +9.74%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '0.590'
$ws.Range("E8").Value = '  +3.64%  '
$ws.Range("D9").Value = '2.432.72'
$ws.Range("E9").Value = '  +5.67%  '
$ws.Range("E10").Value = '  +4.56%  '
$ws.Range("D11").Value = '5.73'
$ws.Range("E11").Value = '  +3.98%  '
$ws.Range("D12").Value = '0.151'
$ws.Range("E12").Value = '  +0.85%  '
$ws.Range("D13").Value = '0.352'
$ws.Range("E13").Value = '  +6.71%  '
$ws.Range("D14").Value = '26.33'
$ws.Range("E14").Value = '  +13.87%  '
$ws.Range("D15").Value = '2.866.88'
$ws.Range("E15").Value = '  +6.07%  '
$ws.Range("D16").Value = '63.133.19'
$ws.Range("E16").Value = '  +6.02%  '
$ws.Range("D17").Value = '0.0000142'
$ws.Range("E17").Value = '  +8.70%  '
$ws.Range("D18").Value = '2.427.97'
$ws.Range("E18").Value = '  +5.65%  '
$ws.Range("E19").Value = '  +8.20%  '
$ws.Range("D20").Value = '340.99'
$ws.Range("E20").Value = '  +10.10%  '
$ws.Range("D21").Value = '4.24'
$ws.Range("E21").Value = '  +5.07%  '
$ws.Range("D22").Value = '6.81'
$ws.Range("E22").Value = '  +4.90%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '65.31'
$ws.Range("E24").Value = '  +3.94%  '
$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").Value = '0.174'
$ws.Range("E25").Value = '  +3.89%  '
$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("B27").Value = 'Fetch.AI'
$ws.Range("C27").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D27").Value = '1.54'
$ws.Range("E27").Value = '  +14.53%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = '8.19'
$ws.Range("E28").Value = '  +6.25%  '
$ws.Range("B29").Value = 'SuiNetwork'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D29").Value = '1.35'
$ws.Range("E29").Value = '  +13.56%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0804'
$ws.Range("E30").Value = '  +12.52%  '
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").Value = '6.69'
$ws.Range("E31").Value = '  +15.75%  '
$ws.Range("D32").Value = '1.83'
$ws.Range("E32").Value = '  +7.21%  '
$ws.Range("B33").Value = 'Monero'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D33").Value = '174.30'
$ws.Range("E33").Value = '  +1.69%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '1.49'
$ws.Range("E34").Value = '  +12.22%  '
$ws.Range("B35").Value = 'PolygonEcosystemToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D35").Value = '0.399'
$ws.Range("E35").Value = '  +6.12%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '18.68'
$ws.Range("E36").Value = '  +6.29%  '
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").Value = '372.43'
$ws.Range("E37").Value = '  +19.12%  '
$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").Value = '4.49'
$ws.Range("E38").Value = '  +13.06%  '
$ws.Range("B39").Value = 'USDe'
$ws.Range("C39").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '1.70'
$ws.Range("E41").Value = '  +13.88%  '
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").Value = '39.91'
$ws.Range("E42").Value = '  +6.67%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '147.72'
$ws.Range("E43").Value = '  +8.55%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").Value = '3.70'
$ws.Range("E44").Value = '  +8.48%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '20.71'
$ws.Range("E45").Value = '  +11.94%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = '0.0959'
$ws.Range("E46").Value = '  +2.16%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '0.593'
$ws.Range("E47").Value = '  +5.10%  '
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").Value = '0.0523'
$ws.Range("E48").Value = '  +7.30%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '0.0225'
$ws.Range("E49").Value = '  +6.56%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '17.88'
$ws.Range("E50").Value = '  +7.57%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").Value = '1.74'
$ws.Range("E51").Value = '  +17.00%  '
